$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 121, pushing the existing
# rows 121-131 down to 123-133 (formatting, e.g. the date style on column D,
# comes along with the insert).
$ws.Range("A121:A122").EntireRow.Insert()

# Fill in the new row 121 with fresh weekly data.
$ws.Range("A121").Value = 3
$ws.Range("B121").Value = "Femacal de La Calera"
$ws.Range("C121").Value = "Coquimbo"
$ws.Range("D121").Value = 44491
$ws.Range("E121").Value = 5
$ws.Range("F121").Value = "Fruta"
$ws.Range("G121").Value = 100101
$ws.Range("H121").Value = "Berries"
$ws.Range("I121").Value = 100112025
$ws.Range("J121").Value = "Frutilla"
$ws.Range("K121").Value = "Sin especificar"
$ws.Range("L121").Value = "Especial"
$ws.Range("M121").Value = 70
$ws.Range("N121").Value = 7000
$ws.Range("O121").Value = 7000
$ws.Range("P121").Value = 7000
$ws.Range("Q121").Value = "$/bandeja 7 kilos"
$ws.Range("R121").Value = "Provincia de Cautín"
$ws.Range("S121").Value = 1000
$ws.Range("T121").Value = 7

# Fill in the new row 122 with fresh weekly data.
$ws.Range("A122").Value = 3
$ws.Range("B122").Value = "Femacal de La Calera"
$ws.Range("C122").Value = "Coquimbo"
$ws.Range("D122").Value = 44491
$ws.Range("E122").Value = 5
$ws.Range("F122").Value = "Fruta"
$ws.Range("G122").Value = 100101
$ws.Range("H122").Value = "Berries"
$ws.Range("I122").Value = 100112025
$ws.Range("J122").Value = "Frutilla"
$ws.Range("K122").Value = "Sin especificar"
$ws.Range("L122").Value = "Segunda"
$ws.Range("M122").Value = 60
$ws.Range("N122").Value = 5000
$ws.Range("O122").Value = 5000
$ws.Range("P122").Value = 5000
$ws.Range("Q122").Value = "$/bandeja 7 kilos"
$ws.Range("R122").Value = "Provincia de Cautín"
$ws.Range("S122").Value = 714
$ws.Range("T122").Value = 7
